$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 8 (2025) metrics with the refreshed figures
$ws.Range("C8").Value = 1108
$ws.Range("D8").Value = 180
$ws.Range("E8").Value = 928
$ws.Range("F8").Value = 7.383100902378999
$ws.Range("G8").Value = 83.75451263537906
$ws.Range("H8").Value = 16.24548736462094
